$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New distribution-label column (AD/AE) for the "best coeff" comparison block ---
$ws.Range("AD19").Value = "0.75_0.1"
$ws.Range("AD20").Value = "0.5_0.1"
$ws.Range("AD21").Value = "0.25_0.1"
$ws.Range("AD22").Value = "uniform"

# --- Section headers ---
$ws.Range("A18").Value = "same poly (best for 8 bit 0.75 0.1) with different input distribution"
$ws.Range("A9").Value = "same poly (best for 16 bit 0.75 0.1) with different bit width"

# --- Annotation note ---
$ws.Range("AE20").Value = "though good than for 0.75 0.1, it is not the best for 0.5 0.1"

# --- Bottom section headers ---
$ws.Range("A26").Value = "with appr mul"
$ws.Range("A24").Value = "all above are precise mul"

# --- Row 19: 8 bit, gaussian 0.75_0.1 ---
$ws.Range("A19").Value = 8
$ws.Range("B19").Value = 4
$ws.Range("C19").Value = 0.25
$ws.Range("D19").Value = 0.875
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.5
$ws.Range("H19").Value = 0.25
$ws.Range("I19").Value = -0.125
$ws.Range("J19").Value = 0.0625
$ws.Range("K19").Value = 0.03125
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 2
$ws.Range("P19").Value = 3
$ws.Range("Q19").Value = 4
$ws.Range("R19").Value = 4
$ws.Range("S19").Value = 4
$ws.Range("T19").Value = 4
$ws.Range("U19").Value = 100000000000000
$ws.Range("V19").Value = 0.14902000000000001
$ws.Range("W19").Value = 0.036343
$ws.Range("X19").Value = 0.006162
$ws.Range("Y19").Value = 0.009827
$ws.Range("Z19").Value = 0.003424
$ws.Range("AA19").Value = 0.002176
$ws.Range("AB19").Value = 0.002176

# --- Row 20: same poly, gaussian 0.5_0.1 ---
$ws.Range("A20").Value = 8
$ws.Range("B20").Value = 4
$ws.Range("C20").Value = 0.25
$ws.Range("D20").Value = 0.875
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.5
$ws.Range("H20").Value = 0.25
$ws.Range("I20").Value = -0.125
$ws.Range("J20").Value = 0.0625
$ws.Range("K20").Value = 0.03125
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 2
$ws.Range("P20").Value = 3
$ws.Range("Q20").Value = 4
$ws.Range("R20").Value = 4
$ws.Range("S20").Value = 4
$ws.Range("T20").Value = 4
$ws.Range("U20").Value = 100000000000000
$ws.Range("V20").Value = 0.091291999999999998
$ws.Range("W20").Value = 0.016688000000000001
$ws.Range("X20").Value = 0.0022929999999999999
$ws.Range("Y20").Value = 0.0023670000000000002
$ws.Range("Z20").Value = 0.001199
$ws.Range("AA20").Value = 0.001165
$ws.Range("AB20").Value = 0.001165

# --- Row 21: same poly, gaussian 0.25_0.1 ---
$ws.Range("A21").Value = 8
$ws.Range("B21").Value = 4
$ws.Range("C21").Value = 0.25
$ws.Range("D21").Value = 0.875
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.5
$ws.Range("H21").Value = 0.25
$ws.Range("I21").Value = -0.125
$ws.Range("J21").Value = 0.0625
$ws.Range("K21").Value = 0.03125
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 2
$ws.Range("P21").Value = 3
$ws.Range("Q21").Value = 4
$ws.Range("R21").Value = 4
$ws.Range("S21").Value = 4
$ws.Range("T21").Value = 4
$ws.Range("U21").Value = 100000000000000
$ws.Range("V21").Value = 0.029173000000000001
$ws.Range("W21").Value = 0.0045560000000000002
$ws.Range("X21").Value = 0.0021870000000000001
$ws.Range("Y21").Value = 0.0021770000000000001
$ws.Range("Z21").Value = 0.0021679999999999998
$ws.Range("AA21").Value = 0.0021679999999999998
$ws.Range("AB21").Value = 0.0021679999999999998

# --- Row 22: same poly, uniform ---
$ws.Range("A22").Value = 8
$ws.Range("B22").Value = 4
$ws.Range("C22").Value = 0.25
$ws.Range("D22").Value = 0.875
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 0.5
$ws.Range("H22").Value = 0.25
$ws.Range("I22").Value = -0.125
$ws.Range("J22").Value = 0.0625
$ws.Range("K22").Value = 0.03125
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 2
$ws.Range("P22").Value = 3
$ws.Range("Q22").Value = 4
$ws.Range("R22").Value = 4
$ws.Range("S22").Value = 4
$ws.Range("T22").Value = 4
$ws.Range("U22").Value = 100000000000000
$ws.Range("V22").Value = 0.074496000000000007
$ws.Range("W22").Value = 0.017056000000000002
$ws.Range("X22").Value = 0.004078
$ws.Range("Y22").Value = 0.0052599999999999999
$ws.Range("Z22").Value = 0.0030509999999999999
$ws.Range("AA22").Value = 0.002624
$ws.Range("AB22").Value = 0.002624

# --- Window / selection bookkeeping to match the saved view state ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("A27").Select()
